$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D1:G235")
# Force text interpretation so numeric-looking codes (e.g. "111") are stored as
# strings rather than being auto-converted to numbers, matching the source data.
$rng.NumberFormat = "@"

$data = New-Object "object[,]" 235,4

$data[0,0] = "codeforiati:group-name"
$data[0,1] = "codeforiati:category-name"
$data[0,2] = "codeforiati:group-code"
$data[0,3] = "codeforiati:category-code"
$data[1,0] = "Educación"
$data[1,1] = "Educación, nivel no especificado"
$data[1,2] = "110"
$data[1,3] = "111"
$data[2,0] = "Educación"
$data[2,1] = "Educación, nivel no especificado"
$data[2,2] = "110"
$data[2,3] = "111"
$data[3,0] = "Educación"
$data[3,1] = "Educación, nivel no especificado"
$data[3,2] = "110"
$data[3,3] = "111"
$data[4,0] = "Educación"
$data[4,1] = "Educación, nivel no especificado"
$data[4,2] = "110"
$data[4,3] = "111"
$data[5,0] = "Educación"
$data[5,1] = "Educación básica"
$data[5,2] = "110"
$data[5,3] = "112"
$data[6,0] = "Educación"
$data[6,1] = "Educación básica"
$data[6,2] = "110"
$data[6,3] = "112"
$data[7,0] = "Educación"
$data[7,1] = "Educación básica"
$data[7,2] = "110"
$data[7,3] = "112"
$data[8,0] = "Educación"
$data[8,1] = "Educación básica"
$data[8,2] = "110"
$data[8,3] = "112"
$data[9,0] = "Educación"
$data[9,1] = "Educación básica"
$data[9,2] = "110"
$data[9,3] = "112"
$data[10,0] = "Educación"
$data[10,1] = "Educación básica"
$data[10,2] = "110"
$data[10,3] = "112"
$data[11,0] = "Educación"
$data[11,1] = "Educación básica"
$data[11,2] = "110"
$data[11,3] = "112"
$data[12,0] = "Educación"
$data[12,1] = "Educación secundaria"
$data[12,2] = "110"
$data[12,3] = "113"
$data[13,0] = "Educación"
$data[13,1] = "Educación secundaria"
$data[13,2] = "110"
$data[13,3] = "113"
$data[14,0] = "Educación"
$data[14,1] = "Educación post-secundaria"
$data[14,2] = "110"
$data[14,3] = "114"
$data[15,0] = "Educación"
$data[15,1] = "Educación post-secundaria"
$data[15,2] = "110"
$data[15,3] = "114"
$data[16,0] = "Salud"
$data[16,1] = "Salud, general"
$data[16,2] = "120"
$data[16,3] = "121"
$data[17,0] = "Salud"
$data[17,1] = "Salud, general"
$data[17,2] = "120"
$data[17,3] = "121"
$data[18,0] = "Salud"
$data[18,1] = "Salud, general"
$data[18,2] = "120"
$data[18,3] = "121"
$data[19,0] = "Salud"
$data[19,1] = "Salud, general"
$data[19,2] = "120"
$data[19,3] = "121"
$data[20,0] = "Salud"
$data[20,1] = "Salud básica"
$data[20,2] = "120"
$data[20,3] = "122"
$data[21,0] = "Salud"
$data[21,1] = "Salud básica"
$data[21,2] = "120"
$data[21,3] = "122"
$data[22,0] = "Salud"
$data[22,1] = "Salud básica"
$data[22,2] = "120"
$data[22,3] = "122"
$data[23,0] = "Salud"
$data[23,1] = "Salud básica"
$data[23,2] = "120"
$data[23,3] = "122"
$data[24,0] = "Salud"
$data[24,1] = "Salud básica"
$data[24,2] = "120"
$data[24,3] = "122"
$data[25,0] = "Salud"
$data[25,1] = "Salud básica"
$data[25,2] = "120"
$data[25,3] = "122"
$data[26,0] = "Salud"
$data[26,1] = "Salud básica"
$data[26,2] = "120"
$data[26,3] = "122"
$data[27,0] = "Salud"
$data[27,1] = "Salud básica"
$data[27,2] = "120"
$data[27,3] = "122"
$data[28,0] = "Salud"
$data[28,1] = "Salud básica"
$data[28,2] = "120"
$data[28,3] = "122"
$data[29,0] = "Salud"
$data[29,1] = "Enfermedades no transmisibles"
$data[29,2] = "120"
$data[29,3] = "123"
$data[30,0] = "Salud"
$data[30,1] = "Enfermedades no transmisibles"
$data[30,2] = "120"
$data[30,3] = "123"
$data[31,0] = "Salud"
$data[31,1] = "Enfermedades no transmisibles"
$data[31,2] = "120"
$data[31,3] = "123"
$data[32,0] = "Salud"
$data[32,1] = "Enfermedades no transmisibles"
$data[32,2] = "120"
$data[32,3] = "123"
$data[33,0] = "Salud"
$data[33,1] = "Enfermedades no transmisibles"
$data[33,2] = "120"
$data[33,3] = "123"
$data[34,0] = "Salud"
$data[34,1] = "Enfermedades no transmisibles"
$data[34,2] = "120"
$data[34,3] = "123"
$data[35,0] = "Programas/políticas sobre población y salud reproductiva"
$data[35,1] = "Programas/políticas sobre población y salud reproductiva"
$data[35,2] = "130"
$data[35,3] = "130"
$data[36,0] = "Programas/políticas sobre población y salud reproductiva"
$data[36,1] = "Programas/políticas sobre población y salud reproductiva"
$data[36,2] = "130"
$data[36,3] = "130"
$data[37,0] = "Programas/políticas sobre población y salud reproductiva"
$data[37,1] = "Programas/políticas sobre población y salud reproductiva"
$data[37,2] = "130"
$data[37,3] = "130"
$data[38,0] = "Programas/políticas sobre población y salud reproductiva"
$data[38,1] = "Programas/políticas sobre población y salud reproductiva"
$data[38,2] = "130"
$data[38,3] = "130"
$data[39,0] = "Programas/políticas sobre población y salud reproductiva"
$data[39,1] = "Programas/políticas sobre población y salud reproductiva"
$data[39,2] = "130"
$data[39,3] = "130"
$data[40,0] = "Abastecimiento de agua y saneamiento"
$data[40,1] = "Abastecimiento de agua y saneamiento"
$data[40,2] = "140"
$data[40,3] = "140"
$data[41,0] = "Abastecimiento de agua y saneamiento"
$data[41,1] = "Abastecimiento de agua y saneamiento"
$data[41,2] = "140"
$data[41,3] = "140"
$data[42,0] = "Abastecimiento de agua y saneamiento"
$data[42,1] = "Abastecimiento de agua y saneamiento"
$data[42,2] = "140"
$data[42,3] = "140"
$data[43,0] = "Abastecimiento de agua y saneamiento"
$data[43,1] = "Abastecimiento de agua y saneamiento"
$data[43,2] = "140"
$data[43,3] = "140"
$data[44,0] = "Abastecimiento de agua y saneamiento"
$data[44,1] = "Abastecimiento de agua y saneamiento"
$data[44,2] = "140"
$data[44,3] = "140"
$data[45,0] = "Abastecimiento de agua y saneamiento"
$data[45,1] = "Abastecimiento de agua y saneamiento"
$data[45,2] = "140"
$data[45,3] = "140"
$data[46,0] = "Abastecimiento de agua y saneamiento"
$data[46,1] = "Abastecimiento de agua y saneamiento"
$data[46,2] = "140"
$data[46,3] = "140"
$data[47,0] = "Abastecimiento de agua y saneamiento"
$data[47,1] = "Abastecimiento de agua y saneamiento"
$data[47,2] = "140"
$data[47,3] = "140"
$data[48,0] = "Abastecimiento de agua y saneamiento"
$data[48,1] = "Abastecimiento de agua y saneamiento"
$data[48,2] = "140"
$data[48,3] = "140"
$data[49,0] = "Abastecimiento de agua y saneamiento"
$data[49,1] = "Abastecimiento de agua y saneamiento"
$data[49,2] = "140"
$data[49,3] = "140"
$data[50,0] = "Abastecimiento de agua y saneamiento"
$data[50,1] = "Abastecimiento de agua y saneamiento"
$data[50,2] = "140"
$data[50,3] = "140"
$data[51,0] = "Gobierno y sociedad civil"
$data[51,1] = "Gobierno y sociedad civil, general"
$data[51,2] = "150"
$data[51,3] = "151"
$data[52,0] = "Gobierno y sociedad civil"
$data[52,1] = "Gobierno y sociedad civil, general"
$data[52,2] = "150"
$data[52,3] = "151"
$data[53,0] = "Gobierno y sociedad civil"
$data[53,1] = "Gobierno y sociedad civil, general"
$data[53,2] = "150"
$data[53,3] = "151"
$data[54,0] = "Gobierno y sociedad civil"
$data[54,1] = "Gobierno y sociedad civil, general"
$data[54,2] = "150"
$data[54,3] = "151"
$data[55,0] = "Gobierno y sociedad civil"
$data[55,1] = "Gobierno y sociedad civil, general"
$data[55,2] = "150"
$data[55,3] = "151"
$data[56,0] = "Gobierno y sociedad civil"
$data[56,1] = "Gobierno y sociedad civil, general"
$data[56,2] = "150"
$data[56,3] = "151"
$data[57,0] = "Gobierno y sociedad civil"
$data[57,1] = "Gobierno y sociedad civil, general"
$data[57,2] = "150"
$data[57,3] = "151"
$data[58,0] = "Gobierno y sociedad civil"
$data[58,1] = "Gobierno y sociedad civil, general"
$data[58,2] = "150"
$data[58,3] = "151"
$data[59,0] = "Gobierno y sociedad civil"
$data[59,1] = "Gobierno y sociedad civil, general"
$data[59,2] = "150"
$data[59,3] = "151"
$data[60,0] = "Gobierno y sociedad civil"
$data[60,1] = "Gobierno y sociedad civil, general"
$data[60,2] = "150"
$data[60,3] = "151"
$data[61,0] = "Gobierno y sociedad civil"
$data[61,1] = "Gobierno y sociedad civil, general"
$data[61,2] = "150"
$data[61,3] = "151"
$data[62,0] = "Gobierno y sociedad civil"
$data[62,1] = "Gobierno y sociedad civil, general"
$data[62,2] = "150"
$data[62,3] = "151"
$data[63,0] = "Gobierno y sociedad civil"
$data[63,1] = "Gobierno y sociedad civil, general"
$data[63,2] = "150"
$data[63,3] = "151"
$data[64,0] = "Gobierno y sociedad civil"
$data[64,1] = "Gobierno y sociedad civil, general"
$data[64,2] = "150"
$data[64,3] = "151"
$data[65,0] = "Gobierno y sociedad civil"
$data[65,1] = "Gobierno y sociedad civil, general"
$data[65,2] = "150"
$data[65,3] = "151"
$data[66,0] = "Gobierno y sociedad civil"
$data[66,1] = "Gobierno y sociedad civil, general"
$data[66,2] = "150"
$data[66,3] = "151"
$data[67,0] = "Gobierno y sociedad civil"
$data[67,1] = "Prevención y resolución de conflictos, paz y seguridad"
$data[67,2] = "150"
$data[67,3] = "152"
$data[68,0] = "Gobierno y sociedad civil"
$data[68,1] = "Prevención y resolución de conflictos, paz y seguridad"
$data[68,2] = "150"
$data[68,3] = "152"
$data[69,0] = "Gobierno y sociedad civil"
$data[69,1] = "Prevención y resolución de conflictos, paz y seguridad"
$data[69,2] = "150"
$data[69,3] = "152"
$data[70,0] = "Gobierno y sociedad civil"
$data[70,1] = "Prevención y resolución de conflictos, paz y seguridad"
$data[70,2] = "150"
$data[70,3] = "152"
$data[71,0] = "Gobierno y sociedad civil"
$data[71,1] = "Prevención y resolución de conflictos, paz y seguridad"
$data[71,2] = "150"
$data[71,3] = "152"
$data[72,0] = "Gobierno y sociedad civil"
$data[72,1] = "Prevención y resolución de conflictos, paz y seguridad"
$data[72,2] = "150"
$data[72,3] = "152"
$data[73,0] = "Otros servicios e infraestructuras sociales"
$data[73,1] = "Otros servicios e infraestructuras sociales"
$data[73,2] = "160"
$data[73,3] = "160"
$data[74,0] = "Otros servicios e infraestructuras sociales"
$data[74,1] = "Otros servicios e infraestructuras sociales"
$data[74,2] = "160"
$data[74,3] = "160"
$data[75,0] = "Otros servicios e infraestructuras sociales"
$data[75,1] = "Otros servicios e infraestructuras sociales"
$data[75,2] = "160"
$data[75,3] = "160"
$data[76,0] = "Otros servicios e infraestructuras sociales"
$data[76,1] = "Otros servicios e infraestructuras sociales"
$data[76,2] = "160"
$data[76,3] = "160"
$data[77,0] = "Otros servicios e infraestructuras sociales"
$data[77,1] = "Otros servicios e infraestructuras sociales"
$data[77,2] = "160"
$data[77,3] = "160"
$data[78,0] = "Otros servicios e infraestructuras sociales"
$data[78,1] = "Otros servicios e infraestructuras sociales"
$data[78,2] = "160"
$data[78,3] = "160"
$data[79,0] = "Otros servicios e infraestructuras sociales"
$data[79,1] = "Otros servicios e infraestructuras sociales"
$data[79,2] = "160"
$data[79,3] = "160"
$data[80,0] = "Otros servicios e infraestructuras sociales"
$data[80,1] = "Otros servicios e infraestructuras sociales"
$data[80,2] = "160"
$data[80,3] = "160"
$data[81,0] = "Otros servicios e infraestructuras sociales"
$data[81,1] = "Otros servicios e infraestructuras sociales"
$data[81,2] = "160"
$data[81,3] = "160"
$data[82,0] = "Otros servicios e infraestructuras sociales"
$data[82,1] = "Otros servicios e infraestructuras sociales"
$data[82,2] = "160"
$data[82,3] = "160"
$data[83,0] = "Otros servicios e infraestructuras sociales"
$data[83,1] = "Otros servicios e infraestructuras sociales"
$data[83,2] = "160"
$data[83,3] = "160"
$data[84,0] = "Transporte y almacenamiento"
$data[84,1] = "Transporte y almacenamiento"
$data[84,2] = "210"
$data[84,3] = "210"
$data[85,0] = "Transporte y almacenamiento"
$data[85,1] = "Transporte y almacenamiento"
$data[85,2] = "210"
$data[85,3] = "210"
$data[86,0] = "Transporte y almacenamiento"
$data[86,1] = "Transporte y almacenamiento"
$data[86,2] = "210"
$data[86,3] = "210"
$data[87,0] = "Transporte y almacenamiento"
$data[87,1] = "Transporte y almacenamiento"
$data[87,2] = "210"
$data[87,3] = "210"
$data[88,0] = "Transporte y almacenamiento"
$data[88,1] = "Transporte y almacenamiento"
$data[88,2] = "210"
$data[88,3] = "210"
$data[89,0] = "Transporte y almacenamiento"
$data[89,1] = "Transporte y almacenamiento"
$data[89,2] = "210"
$data[89,3] = "210"
$data[90,0] = "Transporte y almacenamiento"
$data[90,1] = "Transporte y almacenamiento"
$data[90,2] = "210"
$data[90,3] = "210"
$data[91,0] = "Comunicaciones"
$data[91,1] = "Comunicaciones"
$data[91,2] = "220"
$data[91,3] = "220"
$data[92,0] = "Comunicaciones"
$data[92,1] = "Comunicaciones"
$data[92,2] = "220"
$data[92,3] = "220"
$data[93,0] = "Comunicaciones"
$data[93,1] = "Comunicaciones"
$data[93,2] = "220"
$data[93,3] = "220"
$data[94,0] = "Comunicaciones"
$data[94,1] = "Comunicaciones"
$data[94,2] = "220"
$data[94,3] = "220"
$data[95,0] = "Energía"
$data[95,1] = "Política energética"
$data[95,2] = "230"
$data[95,3] = "231"
$data[96,0] = "Energía"
$data[96,1] = "Política energética"
$data[96,2] = "230"
$data[96,3] = "231"
$data[97,0] = "Energía"
$data[97,1] = "Política energética"
$data[97,2] = "230"
$data[97,3] = "231"
$data[98,0] = "Energía"
$data[98,1] = "Política energética"
$data[98,2] = "230"
$data[98,3] = "231"
$data[99,0] = "Energía"
$data[99,1] = "Generación de energía, fuentes renovables"
$data[99,2] = "230"
$data[99,3] = "232"
$data[100,0] = "Energía"
$data[100,1] = "Generación de energía, fuentes renovables"
$data[100,2] = "230"
$data[100,3] = "232"
$data[101,0] = "Energía"
$data[101,1] = "Generación de energía, fuentes renovables"
$data[101,2] = "230"
$data[101,3] = "232"
$data[102,0] = "Energía"
$data[102,1] = "Generación de energía, fuentes renovables"
$data[102,2] = "230"
$data[102,3] = "232"
$data[103,0] = "Energía"
$data[103,1] = "Generación de energía, fuentes renovables"
$data[103,2] = "230"
$data[103,3] = "232"
$data[104,0] = "Energía"
$data[104,1] = "Generación de energía, fuentes renovables"
$data[104,2] = "230"
$data[104,3] = "232"
$data[105,0] = "Energía"
$data[105,1] = "Generación de energía, fuentes renovables"
$data[105,2] = "230"
$data[105,3] = "232"
$data[106,0] = "Energía"
$data[106,1] = "Generación de energía, fuentes renovables"
$data[106,2] = "230"
$data[106,3] = "232"
$data[107,0] = "Energía"
$data[107,1] = "Generación de energía, fuentes renovables"
$data[107,2] = "230"
$data[107,3] = "232"
$data[108,0] = "Energía"
$data[108,1] = "Generación de energía, fuentes no renovables"
$data[108,2] = "230"
$data[108,3] = "233"
$data[109,0] = "Energía"
$data[109,1] = "Generación de energía, fuentes no renovables"
$data[109,2] = "230"
$data[109,3] = "233"
$data[110,0] = "Energía"
$data[110,1] = "Generación de energía, fuentes no renovables"
$data[110,2] = "230"
$data[110,3] = "233"
$data[111,0] = "Energía"
$data[111,1] = "Generación de energía, fuentes no renovables"
$data[111,2] = "230"
$data[111,3] = "233"
$data[112,0] = "Energía"
$data[112,1] = "Generación de energía, fuentes no renovables"
$data[112,2] = "230"
$data[112,3] = "233"
$data[113,0] = "Energía"
$data[113,1] = "Generación de energía, fuentes no renovables"
$data[113,2] = "230"
$data[113,3] = "233"
$data[114,0] = "Energía"
$data[114,1] = "Centrales de energía híbrida"
$data[114,2] = "230"
$data[114,3] = "234"
$data[115,0] = "Energía"
$data[115,1] = "Centrales de energía nuclear"
$data[115,2] = "230"
$data[115,3] = "235"
$data[116,0] = "Energía"
$data[116,1] = "Distribución de la energía"
$data[116,2] = "230"
$data[116,3] = "236"
$data[117,0] = "Energía"
$data[117,1] = "Distribución de la energía"
$data[117,2] = "230"
$data[117,3] = "236"
$data[118,0] = "Energía"
$data[118,1] = "Distribución de la energía"
$data[118,2] = "230"
$data[118,3] = "236"
$data[119,0] = "Energía"
$data[119,1] = "Distribución de la energía"
$data[119,2] = "230"
$data[119,3] = "236"
$data[120,0] = "Energía"
$data[120,1] = "Distribución de la energía"
$data[120,2] = "230"
$data[120,3] = "236"
$data[121,0] = "Energía"
$data[121,1] = "Distribución de la energía"
$data[121,2] = "230"
$data[121,3] = "236"
$data[122,0] = "Energía"
$data[122,1] = "Distribución de la energía"
$data[122,2] = "230"
$data[122,3] = "236"
$data[123,0] = "Servicios bancarios y financieros"
$data[123,1] = "Servicios bancarios y financieros"
$data[123,2] = "240"
$data[123,3] = "240"
$data[124,0] = "Servicios bancarios y financieros"
$data[124,1] = "Servicios bancarios y financieros"
$data[124,2] = "240"
$data[124,3] = "240"
$data[125,0] = "Servicios bancarios y financieros"
$data[125,1] = "Servicios bancarios y financieros"
$data[125,2] = "240"
$data[125,3] = "240"
$data[126,0] = "Servicios bancarios y financieros"
$data[126,1] = "Servicios bancarios y financieros"
$data[126,2] = "240"
$data[126,3] = "240"
$data[127,0] = "Servicios bancarios y financieros"
$data[127,1] = "Servicios bancarios y financieros"
$data[127,2] = "240"
$data[127,3] = "240"
$data[128,0] = "Servicios bancarios y financieros"
$data[128,1] = "Servicios bancarios y financieros"
$data[128,2] = "240"
$data[128,3] = "240"
$data[129,0] = "Empresas y otros servicios"
$data[129,1] = "Empresas y otros servicios"
$data[129,2] = "250"
$data[129,3] = "250"
$data[130,0] = "Empresas y otros servicios"
$data[130,1] = "Empresas y otros servicios"
$data[130,2] = "250"
$data[130,3] = "250"
$data[131,0] = "Empresas y otros servicios"
$data[131,1] = "Empresas y otros servicios"
$data[131,2] = "250"
$data[131,3] = "250"
$data[132,0] = "Empresas y otros servicios"
$data[132,1] = "Empresas y otros servicios"
$data[132,2] = "250"
$data[132,3] = "250"
$data[133,0] = "Agricultura, Silvicultura, Pesca"
$data[133,1] = "Agricultura"
$data[133,2] = "310"
$data[133,3] = "311"
$data[134,0] = "Agricultura, Silvicultura, Pesca"
$data[134,1] = "Agricultura"
$data[134,2] = "310"
$data[134,3] = "311"
$data[135,0] = "Agricultura, Silvicultura, Pesca"
$data[135,1] = "Agricultura"
$data[135,2] = "310"
$data[135,3] = "311"
$data[136,0] = "Agricultura, Silvicultura, Pesca"
$data[136,1] = "Agricultura"
$data[136,2] = "310"
$data[136,3] = "311"
$data[137,0] = "Agricultura, Silvicultura, Pesca"
$data[137,1] = "Agricultura"
$data[137,2] = "310"
$data[137,3] = "311"
$data[138,0] = "Agricultura, Silvicultura, Pesca"
$data[138,1] = "Agricultura"
$data[138,2] = "310"
$data[138,3] = "311"
$data[139,0] = "Agricultura, Silvicultura, Pesca"
$data[139,1] = "Agricultura"
$data[139,2] = "310"
$data[139,3] = "311"
$data[140,0] = "Agricultura, Silvicultura, Pesca"
$data[140,1] = "Agricultura"
$data[140,2] = "310"
$data[140,3] = "311"
$data[141,0] = "Agricultura, Silvicultura, Pesca"
$data[141,1] = "Agricultura"
$data[141,2] = "310"
$data[141,3] = "311"
$data[142,0] = "Agricultura, Silvicultura, Pesca"
$data[142,1] = "Agricultura"
$data[142,2] = "310"
$data[142,3] = "311"
$data[143,0] = "Agricultura, Silvicultura, Pesca"
$data[143,1] = "Agricultura"
$data[143,2] = "310"
$data[143,3] = "311"
$data[144,0] = "Agricultura, Silvicultura, Pesca"
$data[144,1] = "Agricultura"
$data[144,2] = "310"
$data[144,3] = "311"
$data[145,0] = "Agricultura, Silvicultura, Pesca"
$data[145,1] = "Agricultura"
$data[145,2] = "310"
$data[145,3] = "311"
$data[146,0] = "Agricultura, Silvicultura, Pesca"
$data[146,1] = "Agricultura"
$data[146,2] = "310"
$data[146,3] = "311"
$data[147,0] = "Agricultura, Silvicultura, Pesca"
$data[147,1] = "Agricultura"
$data[147,2] = "310"
$data[147,3] = "311"
$data[148,0] = "Agricultura, Silvicultura, Pesca"
$data[148,1] = "Agricultura"
$data[148,2] = "310"
$data[148,3] = "311"
$data[149,0] = "Agricultura, Silvicultura, Pesca"
$data[149,1] = "Agricultura"
$data[149,2] = "310"
$data[149,3] = "311"
$data[150,0] = "Agricultura, Silvicultura, Pesca"
$data[150,1] = "Agricultura"
$data[150,2] = "310"
$data[150,3] = "311"
$data[151,0] = "Agricultura, Silvicultura, Pesca"
$data[151,1] = "Silvicultura"
$data[151,2] = "310"
$data[151,3] = "312"
$data[152,0] = "Agricultura, Silvicultura, Pesca"
$data[152,1] = "Silvicultura"
$data[152,2] = "310"
$data[152,3] = "312"
$data[153,0] = "Agricultura, Silvicultura, Pesca"
$data[153,1] = "Silvicultura"
$data[153,2] = "310"
$data[153,3] = "312"
$data[154,0] = "Agricultura, Silvicultura, Pesca"
$data[154,1] = "Silvicultura"
$data[154,2] = "310"
$data[154,3] = "312"
$data[155,0] = "Agricultura, Silvicultura, Pesca"
$data[155,1] = "Silvicultura"
$data[155,2] = "310"
$data[155,3] = "312"
$data[156,0] = "Agricultura, Silvicultura, Pesca"
$data[156,1] = "Silvicultura"
$data[156,2] = "310"
$data[156,3] = "312"
$data[157,0] = "Agricultura, Silvicultura, Pesca"
$data[157,1] = "Pesca"
$data[157,2] = "310"
$data[157,3] = "313"
$data[158,0] = "Agricultura, Silvicultura, Pesca"
$data[158,1] = "Pesca"
$data[158,2] = "310"
$data[158,3] = "313"
$data[159,0] = "Agricultura, Silvicultura, Pesca"
$data[159,1] = "Pesca"
$data[159,2] = "310"
$data[159,3] = "313"
$data[160,0] = "Agricultura, Silvicultura, Pesca"
$data[160,1] = "Pesca"
$data[160,2] = "310"
$data[160,3] = "313"
$data[161,0] = "Agricultura, Silvicultura, Pesca"
$data[161,1] = "Pesca"
$data[161,2] = "310"
$data[161,3] = "313"
$data[162,0] = "Industria, extractivas, construcción"
$data[162,1] = "Industria"
$data[162,2] = "320"
$data[162,3] = "321"
$data[163,0] = "Industria, extractivas, construcción"
$data[163,1] = "Industria"
$data[163,2] = "320"
$data[163,3] = "321"
$data[164,0] = "Industria, extractivas, construcción"
$data[164,1] = "Industria"
$data[164,2] = "320"
$data[164,3] = "321"
$data[165,0] = "Industria, extractivas, construcción"
$data[165,1] = "Industria"
$data[165,2] = "320"
$data[165,3] = "321"
$data[166,0] = "Industria, extractivas, construcción"
$data[166,1] = "Industria"
$data[166,2] = "320"
$data[166,3] = "321"
$data[167,0] = "Industria, extractivas, construcción"
$data[167,1] = "Industria"
$data[167,2] = "320"
$data[167,3] = "321"
$data[168,0] = "Industria, extractivas, construcción"
$data[168,1] = "Industria"
$data[168,2] = "320"
$data[168,3] = "321"
$data[169,0] = "Industria, extractivas, construcción"
$data[169,1] = "Industria"
$data[169,2] = "320"
$data[169,3] = "321"
$data[170,0] = "Industria, extractivas, construcción"
$data[170,1] = "Industria"
$data[170,2] = "320"
$data[170,3] = "321"
$data[171,0] = "Industria, extractivas, construcción"
$data[171,1] = "Industria"
$data[171,2] = "320"
$data[171,3] = "321"
$data[172,0] = "Industria, extractivas, construcción"
$data[172,1] = "Industria"
$data[172,2] = "320"
$data[172,3] = "321"
$data[173,0] = "Industria, extractivas, construcción"
$data[173,1] = "Industria"
$data[173,2] = "320"
$data[173,3] = "321"
$data[174,0] = "Industria, extractivas, construcción"
$data[174,1] = "Industria"
$data[174,2] = "320"
$data[174,3] = "321"
$data[175,0] = "Industria, extractivas, construcción"
$data[175,1] = "Industria"
$data[175,2] = "320"
$data[175,3] = "321"
$data[176,0] = "Industria, extractivas, construcción"
$data[176,1] = "Industria"
$data[176,2] = "320"
$data[176,3] = "321"
$data[177,0] = "Industria, extractivas, construcción"
$data[177,1] = "Industria"
$data[177,2] = "320"
$data[177,3] = "321"
$data[178,0] = "Industria, extractivas, construcción"
$data[178,1] = "Industria"
$data[178,2] = "320"
$data[178,3] = "321"
$data[179,0] = "Industria, extractivas, construcción"
$data[179,1] = "Industria"
$data[179,2] = "320"
$data[179,3] = "321"
$data[180,0] = "Industria, extractivas, construcción"
$data[180,1] = "Industria"
$data[180,2] = "320"
$data[180,3] = "321"
$data[181,0] = "Industria, extractivas, construcción"
$data[181,1] = "Industrias extractivas"
$data[181,2] = "320"
$data[181,3] = "322"
$data[182,0] = "Industria, extractivas, construcción"
$data[182,1] = "Industrias extractivas"
$data[182,2] = "320"
$data[182,3] = "322"
$data[183,0] = "Industria, extractivas, construcción"
$data[183,1] = "Industrias extractivas"
$data[183,2] = "320"
$data[183,3] = "322"
$data[184,0] = "Industria, extractivas, construcción"
$data[184,1] = "Industrias extractivas"
$data[184,2] = "320"
$data[184,3] = "322"
$data[185,0] = "Industria, extractivas, construcción"
$data[185,1] = "Industrias extractivas"
$data[185,2] = "320"
$data[185,3] = "322"
$data[186,0] = "Industria, extractivas, construcción"
$data[186,1] = "Industrias extractivas"
$data[186,2] = "320"
$data[186,3] = "322"
$data[187,0] = "Industria, extractivas, construcción"
$data[187,1] = "Industrias extractivas"
$data[187,2] = "320"
$data[187,3] = "322"
$data[188,0] = "Industria, extractivas, construcción"
$data[188,1] = "Industrias extractivas"
$data[188,2] = "320"
$data[188,3] = "322"
$data[189,0] = "Industria, extractivas, construcción"
$data[189,1] = "Industrias extractivas"
$data[189,2] = "320"
$data[189,3] = "322"
$data[190,0] = "Industria, extractivas, construcción"
$data[190,1] = "Industrias extractivas"
$data[190,2] = "320"
$data[190,3] = "322"
$data[191,0] = "Industria, extractivas, construcción"
$data[191,1] = "Construcción"
$data[191,2] = "320"
$data[191,3] = "323"
$data[192,0] = "Política y regulación comercial"
$data[192,1] = "Política y regulación comercial"
$data[192,2] = "331"
$data[192,3] = "331"
$data[193,0] = "Política y regulación comercial"
$data[193,1] = "Política y regulación comercial"
$data[193,2] = "331"
$data[193,3] = "331"
$data[194,0] = "Política y regulación comercial"
$data[194,1] = "Política y regulación comercial"
$data[194,2] = "331"
$data[194,3] = "331"
$data[195,0] = "Política y regulación comercial"
$data[195,1] = "Política y regulación comercial"
$data[195,2] = "331"
$data[195,3] = "331"
$data[196,0] = "Política y regulación comercial"
$data[196,1] = "Política y regulación comercial"
$data[196,2] = "331"
$data[196,3] = "331"
$data[197,0] = "Política y regulación comercial"
$data[197,1] = "Política y regulación comercial"
$data[197,2] = "331"
$data[197,3] = "331"
$data[198,0] = "Turismo"
$data[198,1] = "Turismo"
$data[198,2] = "332"
$data[198,3] = "332"
$data[199,0] = "Protección general medio ambiente"
$data[199,1] = "Protección general medio ambiente"
$data[199,2] = "410"
$data[199,3] = "410"
$data[200,0] = "Protección general medio ambiente"
$data[200,1] = "Protección general medio ambiente"
$data[200,2] = "410"
$data[200,3] = "410"
$data[201,0] = "Protección general medio ambiente"
$data[201,1] = "Protección general medio ambiente"
$data[201,2] = "410"
$data[201,3] = "410"
$data[202,0] = "Protección general medio ambiente"
$data[202,1] = "Protección general medio ambiente"
$data[202,2] = "410"
$data[202,3] = "410"
$data[203,0] = "Protección general medio ambiente"
$data[203,1] = "Protección general medio ambiente"
$data[203,2] = "410"
$data[203,3] = "410"
$data[204,0] = "Protección general medio ambiente"
$data[204,1] = "Protección general medio ambiente"
$data[204,2] = "410"
$data[204,3] = "410"
$data[205,0] = "Otras acciones multisectoriales"
$data[205,1] = "Otras acciones multisectoriales"
$data[205,2] = "430"
$data[205,3] = "430"
$data[206,0] = "Otras acciones multisectoriales"
$data[206,1] = "Otras acciones multisectoriales"
$data[206,2] = "430"
$data[206,3] = "430"
$data[207,0] = "Otras acciones multisectoriales"
$data[207,1] = "Otras acciones multisectoriales"
$data[207,2] = "430"
$data[207,3] = "430"
$data[208,0] = "Otras acciones multisectoriales"
$data[208,1] = "Otras acciones multisectoriales"
$data[208,2] = "430"
$data[208,3] = "430"
$data[209,0] = "Otras acciones multisectoriales"
$data[209,1] = "Otras acciones multisectoriales"
$data[209,2] = "430"
$data[209,3] = "430"
$data[210,0] = "Otras acciones multisectoriales"
$data[210,1] = "Otras acciones multisectoriales"
$data[210,2] = "430"
$data[210,3] = "430"
$data[211,0] = "Otras acciones multisectoriales"
$data[211,1] = "Otras acciones multisectoriales"
$data[211,2] = "430"
$data[211,3] = "430"
$data[212,0] = "Otras acciones multisectoriales"
$data[212,1] = "Otras acciones multisectoriales"
$data[212,2] = "430"
$data[212,3] = "430"
$data[213,0] = "Otras acciones multisectoriales"
$data[213,1] = "Otras acciones multisectoriales"
$data[213,2] = "430"
$data[213,3] = "430"
$data[214,0] = "Otras acciones multisectoriales"
$data[214,1] = "Otras acciones multisectoriales"
$data[214,2] = "430"
$data[214,3] = "430"
$data[215,0] = "Apoyo presupuestario general"
$data[215,1] = "Apoyo presupuestario general"
$data[215,2] = "510"
$data[215,3] = "510"
$data[216,0] = "Ayuda alimentaria para el desarrollo"
$data[216,1] = "Ayuda alimentaria para el desarrollo"
$data[216,2] = "520"
$data[216,3] = "520"
$data[217,0] = "Otras ayudas en forma de suministro de bienes"
$data[217,1] = "Otras ayudas en forma de suministro de bienes"
$data[217,2] = "530"
$data[217,3] = "530"
$data[218,0] = "Otras ayudas en forma de suministro de bienes"
$data[218,1] = "Otras ayudas en forma de suministro de bienes"
$data[218,2] = "530"
$data[218,3] = "530"
$data[219,0] = "Actividades relacionadas con la deuda"
$data[219,1] = "Actividades relacionadas con la deuda"
$data[219,2] = "600"
$data[219,3] = "600"
$data[220,0] = "Actividades relacionadas con la deuda"
$data[220,1] = "Actividades relacionadas con la deuda"
$data[220,2] = "600"
$data[220,3] = "600"
$data[221,0] = "Actividades relacionadas con la deuda"
$data[221,1] = "Actividades relacionadas con la deuda"
$data[221,2] = "600"
$data[221,3] = "600"
$data[222,0] = "Actividades relacionadas con la deuda"
$data[222,1] = "Actividades relacionadas con la deuda"
$data[222,2] = "600"
$data[222,3] = "600"
$data[223,0] = "Actividades relacionadas con la deuda"
$data[223,1] = "Actividades relacionadas con la deuda"
$data[223,2] = "600"
$data[223,3] = "600"
$data[224,0] = "Actividades relacionadas con la deuda"
$data[224,1] = "Actividades relacionadas con la deuda"
$data[224,2] = "600"
$data[224,3] = "600"
$data[225,0] = "Actividades relacionadas con la deuda"
$data[225,1] = "Actividades relacionadas con la deuda"
$data[225,2] = "600"
$data[225,3] = "600"
$data[226,0] = "Ayuda de emergencia"
$data[226,1] = "Ayuda de emergencia"
$data[226,2] = "720"
$data[226,3] = "720"
$data[227,0] = "Ayuda de emergencia"
$data[227,1] = "Ayuda de emergencia"
$data[227,2] = "720"
$data[227,3] = "720"
$data[228,0] = "Ayuda de emergencia"
$data[228,1] = "Ayuda de emergencia"
$data[228,2] = "720"
$data[228,3] = "720"
$data[229,0] = "Ayuda a la reconstrucción y a la rehabilitación"
$data[229,1] = "Ayuda a la reconstrucción y a la rehabilitación"
$data[229,2] = "730"
$data[229,3] = "730"
$data[230,0] = "Prevención de desastres"
$data[230,1] = "Prevención de desastres"
$data[230,2] = "740"
$data[230,3] = "740"
$data[231,0] = "Costes administrativos donantes"
$data[231,1] = "Costes administrativos donantes"
$data[231,2] = "910"
$data[231,3] = "910"
$data[232,0] = "Ayuda a refugiados en el país donante"
$data[232,1] = "Ayuda a refugiados en el país donante"
$data[232,2] = "930"
$data[232,3] = "930"
$data[233,0] = "Sin especificación / no clasificados"
$data[233,1] = "Sin especificación / no clasificados"
$data[233,2] = "998"
$data[233,3] = "998"
$data[234,0] = "Sin especificación / no clasificados"
$data[234,1] = "Sin especificación / no clasificados"
$data[234,2] = "998"
$data[234,3] = "998"

$rng.Value = $data

# Revert the temporary text number-format back to the workbook default so the
# cells end up styled exactly as they were (format General), only the cached
# "treat as text" flag is dropped, not the stored string values.
$rng.Style = "Normal"
